$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted as row 5 (pushing the
# existing rows 5 and 6 down to rows 6 and 7). Insert a blank row at 5 so
# the two trailing records keep their original values/positions relative
# to each other, just shifted down by one.
$ws.Rows.Item(5).Insert()

# Seed the new row 5 from the row right below it (the old row 5, now at
# row 6) so every static/shared column (market, region, product, unit,
# origin, classification, ...) and the date cell's number format come
# along for free.
$ws.Range("A6:R6").Copy() | Out-Null
$ws.Range("A5:R5").PasteSpecial() | Out-Null

# Overwrite just the values that differ for this new weekly record.
$ws.Cells.Item(5, 4).Value = 45205    # Fecha
$ws.Cells.Item(5, 10).Value = 200     # Volumen
$ws.Cells.Item(5, 11).Value = 11000   # Precio mínimo
$ws.Cells.Item(5, 12).Value = 12000   # Precio máximo
$ws.Cells.Item(5, 13).Value = 11500   # Precio promedio ponderado
$ws.Cells.Item(5, 16).Value = 639     # Precio $/Kg

$excel.CutCopyMode = $false
